$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140 (shifts existing rows 140-204 down to 141-205)
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new record
$ws.Range("A140").Value = 9
$ws.Range("B140").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C140").Value = "Metropolitana"
$ws.Range("D140").Value = 44488
$ws.Range("E140").Value = 13
$ws.Range("F140").Value = 100112052
$ws.Range("G140").Value = "Albahaca"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 79
$ws.Range("K140").Value = 7000
$ws.Range("L140").Value = 7000
$ws.Range("M140").Value = 7000
$ws.Range("N140").Value = "`$/docena de matas"
$ws.Range("O140").Value = "Provincia de Chacabuco"
$ws.Range("P140").Value = 1167
$ws.Range("Q140").Value = 6
$ws.Range("R140").Value = "Hortaliza"
